# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 6
    3  = 4
    4  = 7
    5  = 3
    6  = 4
    7  = 7
    8  = 5
    9  = 9
    10 = 4
    11 = 5
    12 = 6
    13 = 4
    14 = 3
    15 = 10
    16 = 9
    17 = 6
    18 = 5
    19 = 8
    20 = 4
    21 = 4
    22 = 7
    23 = 4
    24 = 5
    25 = 2
    26 = 9
    27 = 5
    28 = 6
    29 = 5
    30 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
